# Updated Testcases with Clearing test cases
# Adds a new manual test case (row 41) covering: "Clearing of Bid by
# student who placed $10" and moves the sheet view/selection down to it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Manual Testing")

[void]$ws.Activate()

# --- New test case row (row 41) -----------------------------------------
$ws.Cells.Item(41, 1).Value = 40
$ws.Cells.Item(41, 2).Value = "Clearing of Bid by student who placed `$10"
$ws.Cells.Item(41, 3).Value = "Bid for a course: round 1"
$ws.Cells.Item(41, 4).Value = "Admin logged in and presses Clear Round 1"
$ws.Cells.Item(41, 5).Value = "Student bids for IS100 with `$10 and admin clears it after"
$ws.Cells.Item(41, 6).Value = "Student: amy.ng.2009`nCourse: IS100`nSection: S1`nAmount: `$10"
$ws.Cells.Item(41, 7).Value = "amy.ng.2009 in student table will have `$190 instead of `$200"
$ws.Cells.Item(41, 8).Value = "amy.ng.2009 in student table will have `$190 instead of `$200"
$ws.Cells.Item(41, 9).Value = "Pass"

$ws.Rows.Item(41).RowHeight = 58

# --- Move the viewport / selection down to the new row ------------------
$excel.ActiveWindow.ScrollRow = 40
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("A42:F42").Select()
